$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

for ($r = 15; $r -ge 7; $r--) {
    $srcRow = $r
    $dstRow = $r + 2
    $ws.Range("A$srcRow`:D$srcRow").Copy($ws.Range("A$dstRow"))
}
Write-Host "done shifting"
